$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Colour constants (OLE/COM colour is 0xBBGGRR, i.e. reversed RGB)
# ---------------------------------------------------------------------------
$colGreen     = 5287936   # FF00B050 -> "recu"     fill (fillId 2 / cellXf 4)
$colYellow    = 65535     # FFFFFF00 -> "expedie"  fill (fillId 3 / cellXf 5)
$colRed       = 192       # FFC00000 -> "paye"     fill (fillId 4 / cellXf 6)
$colLightBlue = 15773696  # FF00B0F0 -> "point relais" fill (new fillId 5 / cellXf 8)

# ---------------------------------------------------------------------------
# 1) Three new purchase rows (55-57)
# ---------------------------------------------------------------------------

# Row 55 : Le canon de l'Etoile de la Mort - 75246, bought via Vinted, seller "nestarit"
$ws.Range("A55").Value = 2023
$ws.Range("B55").Value = "https://www.idealo.fr/prix/6799851/lego-star-wars-le-canon-de-l-etoile-de-la-mort-75246.html"
$ws.Range("D55").Value = 75246
$ws.Range("E55").Value = "Vinted"
$ws.Range("G55").Value = 25
$ws.Range("H55").Interior.Color = $colRed
$ws.Range("L55").Value = "non"
$ws.Range("M55").Value = 1
$ws.Range("N55").Value = "nestarit"
$ws.Range("O55").Value = "King jouet"

# Row 56 : Microfighter Y-Wing de la Resistance - 75263, bought via Vinted, seller "nestarit"
$ws.Range("A56").Value = 2023
$ws.Range("B56").Value = "https://www.idealo.fr/prix/6982375/lego-star-wars-microfighter-y-wing-de-la-resistance-75263.html"
$ws.Range("D56").Value = 75263
$ws.Range("E56").Value = "Vinted"
$ws.Range("G56").Value = 11.13
$ws.Range("H56").Interior.Color = $colRed
$ws.Range("L56").Value = "non"
$ws.Range("M56").Value = 1
$ws.Range("N56").Value = "nestarit"
$ws.Range("O56").Value = "King jouet"

# Row 57 : LEGO Microfighter Navette de Kylo Ren(TM) - 75264, bought on Amazon, 4 exemplaires
$ws.Range("A57").Value = 2023
$ws.Range("B57").Value = "https://www.idealo.fr/prix/6982410/lego-star-wars-microfighter-navette-de-kylo-ren-75264.html"
$ws.Range("C57").Value = "LEGO Microfighter Navette de Kylo Ren" + [char]0x2122
$ws.Range("D57").Value = 75264
$ws.Range("E57").Value = "Amazon"
$ws.Range("F57").Value = 29.83
$ws.Range("G57").Value = 119.32
$ws.Range("H57").Interior.Color = $colRed
$ws.Range("L57").Value = "non"
$ws.Range("M57").Value = 4
$ws.Range("N57").Value = "B07W8XYZ2X"
$ws.Range("O57").Value = "Chez moi"

# ---------------------------------------------------------------------------
# 2) Order-status pipeline progresses: a new "point relais" stage is inserted
#    ahead of the existing "expedie" / "paye" entries, shifting them down.
# ---------------------------------------------------------------------------

# Row 3 : was "expedie" (yellow) -> now "point relais" (light blue)
$ws.Range("P3").Value = "point relais"
$ws.Range("Q3").Interior.Color = $colLightBlue

# Row 4 : was "paye" (red) -> now "expedie" (yellow)
$ws.Range("P4").Value = "expédié"
$ws.Range("Q4").Interior.Color = $colYellow

# Row 5 : newly tracked -> "paye" (red)
$ws.Range("P5").Value = "payé"
$ws.Range("Q5").Interior.Color = $colRed

# ---------------------------------------------------------------------------
# 3) Several previously-yellow "Stoke" cells have been received: turn green.
# ---------------------------------------------------------------------------
$greenRows = 26,27,28,29,30,31,32,33,37,38,39
foreach ($r in $greenRows) {
    $ws.Range("H$r").Interior.Color = $colGreen
}

# Row 41 moves all the way to "point relais" (light blue)
$ws.Range("H41").Interior.Color = $colLightBlue

# Rows 50-51 move from "paye" (red) to "expedie" (yellow)
$ws.Range("H50").Interior.Color = $colYellow
$ws.Range("H51").Interior.Color = $colYellow
